$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the existing header style/format from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add data values for the new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
